$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl17"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.244729
$ws.Range("H2").Value = 0.734187
$ws.Range("I2").Value = 0.4052712693903822
$ws.Range("J2").Value = 0.4052712693903822
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.062395
$ws.Range("N2").Value = 0.187185
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.015269865955
$ws.Range("R2").Value = 0.137428793595
$ws.Range("S2").Value = 0.4052712693903822
$ws.Range("T2").Value = 0.4052712693903822

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl17"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.149358
$ws.Range("H3").Value = 0.448074
$ws.Range("I3").Value = 0.2473368757017301
$ws.Range("J3").Value = 0.24733687570173
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.062395
$ws.Range("N3").Value = 0.187185
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.00931919241
$ws.Range("R3").Value = 0.08387273169
$ws.Range("S3").Value = 0.2473368757017301
$ws.Range("T3").Value = 0.24733687570173

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Ccl17"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1021816666666667
$ws.Range("H4").Value = 0.306545
$ws.Range("I4").Value = 0.1692128589518402
$ws.Range("J4").Value = 0.1692128589518402
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.062395
$ws.Range("N4").Value = 0.187185
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.006375625091666667
$ws.Range("R4").Value = 0.057380625825
$ws.Range("S4").Value = 0.1692128589518402
$ws.Range("T4").Value = 0.1692128589518402

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl17"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.107596
$ws.Range("H5").Value = 0.322788
$ws.Range("I5").Value = 0.1781789959560476
$ws.Range("J5").Value = 0.1781789959560476
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.062395
$ws.Range("N5").Value = 0.187185
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.006713452420000001
$ws.Range("R5").Value = 0.06042107178
$ws.Range("S5").Value = 0.1781789959560476
$ws.Range("T5").Value = 0.1781789959560476
